# "Generate Report for Archive"
# - Flip the localization status text from "Ready for handoff" to
#   "In Translation" everywhere it appears (Overview zh-cn/de-de status
#   columns plus the per-language Status column on each language sheet).
# - Narrow the now-shorter "Status" columns (Overview!E:F and the
#   "Status" column on the zh-cn / de-de sheets) to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C1").EntireColumn.ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"
$dede.Range("C1").EntireColumn.ColumnWidth = 12.5
